$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Size (Cities)" column between "Data Set" and "Time (s)":
# shift the existing "Time (s)" header from B10 to C10, then set B10.
$ws.Range("C10").Value = "Time (s)"
$ws.Range("B10").Value = "Size (Cities)"

# Fill in the new numeric data for the small results table
$ws.Range("B11").Value = 52
$ws.Range("C11").Value = 353

$ws.Range("B12").Value = 734
$ws.Range("C12").Value = 7079

$ws.Range("B13").Value = 16862

# Column B now holds its own (narrower) data, separate from column A;
# split it out of the shared A:B column-width group and size it to fit.
$ws.Columns.Item(2).ColumnWidth = 10.7109375

# Restore the active selection to H10
$ws.Range("H10").Select()
